$d = $word.ActiveDocument
$d.Content.Find.Execute("99-77=22", $true, $true, $false, $false, $false, $true, 1, $false, "65-33=32", 2) | Out-Null
$d.Content.Find.Execute("60-36=24", $true, $true, $false, $false, $false, $true, 1, $false, "68+11=79", 2) | Out-Null
$d.Content.Find.Execute("54-50=4", $true, $true, $false, $false, $false, $true, 1, $false, "59+34=93", 2) | Out-Null
$d.Content.Find.Execute("43-7=36", $true, $true, $false, $false, $false, $true, 1, $false, "92-5=87", 2) | Out-Null
$d.Content.Find.Execute("65+6=71", $true, $true, $false, $false, $false, $true, 1, $false, "53+6=59", 2) | Out-Null
$d.Content.Find.Execute("98-18=80", $true, $true, $false, $false, $false, $true, 1, $false, "62-9=53", 2) | Out-Null
$d.Content.Find.Execute("20+31=51", $true, $true, $false, $false, $false, $true, 1, $false, "64+33=97", 2) | Out-Null
$d.Content.Find.Execute("81-10=71", $true, $true, $false, $false, $false, $true, 1, $false, "66-16=50", 2) | Out-Null
$d.Content.Find.Execute("91-8=83", $true, $true, $false, $false, $false, $true, 1, $false, "29+47=76", 2) | Out-Null
$d.Content.Find.Execute("8+68=76", $true, $true, $false, $false, $false, $true, 1, $false, "30+53=83", 2) | Out-Null
$d.Content.Find.Execute("48+17=65", $true, $true, $false, $false, $false, $true, 1, $false, "81-65=16", 2) | Out-Null
$d.Content.Find.Execute("86-36=50", $true, $true, $false, $false, $false, $true, 1, $false, "65-6=59", 2) | Out-Null
$d.Content.Find.Execute("61-61=0", $true, $true, $false, $false, $false, $true, 1, $false, "92-85=7", 2) | Out-Null
$d.Content.Find.Execute("30-21=9", $true, $true, $false, $false, $false, $true, 1, $false, "36+2=38", 2) | Out-Null
$d.Content.Find.Execute("73-37=36", $true, $true, $false, $false, $false, $true, 1, $false, "16+47=63", 2) | Out-Null
$d.Content.Find.Execute("29+63=92", $true, $true, $false, $false, $false, $true, 1, $false, "17+9=26", 2) | Out-Null
$d.Content.Find.Execute("39+16=55", $true, $true, $false, $false, $false, $true, 1, $false, "99-3=96", 2) | Out-Null
$d.Content.Find.Execute("66-14=52", $true, $true, $false, $false, $false, $true, 1, $false, "52-30=22", 2) | Out-Null
$d.Content.Find.Execute("4+94=98", $true, $true, $false, $false, $false, $true, 1, $false, "87-41=46", 2) | Out-Null
$d.Content.Find.Execute("32+57=89", $true, $true, $false, $false, $false, $true, 1, $false, "21+59=80", 2) | Out-Null
$d.Content.Find.Execute("17+10=27", $true, $true, $false, $false, $false, $true, 1, $false, "33+56=89", 2) | Out-Null
$d.Content.Find.Execute("13+48=61", $true, $true, $false, $false, $false, $true, 1, $false, "78-6=72", 2) | Out-Null
$d.Content.Find.Execute("69-62=7", $true, $true, $false, $false, $false, $true, 1, $false, "36-27=9", 2) | Out-Null
$d.Content.Find.Execute("74-62=12", $true, $true, $false, $false, $false, $true, 1, $false, "44+36=80", 2) | Out-Null
$d.Content.Find.Execute("16+77=93", $true, $true, $false, $false, $false, $true, 1, $false, "43+52=95", 2) | Out-Null
$d.Content.Find.Execute("11+11=22", $true, $true, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
$d.Content.Find.Execute("64+23=87", $true, $true, $false, $false, $false, $true, 1, $false, "56+11=67", 2) | Out-Null
$d.Content.Find.Execute("44-16=28", $true, $true, $false, $false, $false, $true, 1, $false, "40+15=55", 2) | Out-Null
$d.Content.Find.Execute("75-28=47", $true, $true, $false, $false, $false, $true, 1, $false, "57+11=68", 2) | Out-Null
$d.Content.Find.Execute("76+2=78", $true, $true, $false, $false, $false, $true, 1, $false, "80-43=37", 2) | Out-Null
$d.Content.Find.Execute("81-23=58", $true, $true, $false, $false, $false, $true, 1, $false, "0+21=21", 2) | Out-Null
$d.Content.Find.Execute("71-44=27", $true, $true, $false, $false, $false, $true, 1, $false, "67+31=98", 2) | Out-Null
$d.Content.Find.Execute("33+59=92", $true, $true, $false, $false, $false, $true, 1, $false, "35-22=13", 2) | Out-Null
$d.Content.Find.Execute("79-26=53", $true, $true, $false, $false, $false, $true, 1, $false, "89-16=73", 2) | Out-Null
$d.Content.Find.Execute("16+67=83", $true, $true, $false, $false, $false, $true, 1, $false, "66-56=10", 2) | Out-Null
$d.Content.Find.Execute("43+21=64", $true, $true, $false, $false, $false, $true, 1, $false, "34+2=36", 2) | Out-Null
$d.Content.Find.Execute("18+31=49", $true, $true, $false, $false, $false, $true, 1, $false, "6+19=25", 2) | Out-Null
$d.Content.Find.Execute("5+3=8", $true, $true, $false, $false, $false, $true, 1, $false, "3+90=93", 2) | Out-Null
$d.Content.Find.Execute("13+2=15", $true, $true, $false, $false, $false, $true, 1, $false, "48+2=50", 2) | Out-Null
$d.Content.Find.Execute("50-0=50", $true, $true, $false, $false, $false, $true, 1, $false, "5+59=64", 2) | Out-Null
$d.Content.Find.Execute("66-63=3", $true, $true, $false, $false, $false, $true, 1, $false, "68+9=77", 2) | Out-Null
$d.Content.Find.Execute("12+49=61", $true, $true, $false, $false, $false, $true, 1, $false, "54-32=22", 2) | Out-Null
$d.Content.Find.Execute("48+50=98", $true, $true, $false, $false, $false, $true, 1, $false, "36+34=70", 2) | Out-Null
$d.Content.Find.Execute("70-3=67", $true, $true, $false, $false, $false, $true, 1, $false, "33+31=64", 2) | Out-Null
$d.Content.Find.Execute("94-53=41", $true, $true, $false, $false, $false, $true, 1, $false, "3+90=93", 2) | Out-Null
$d.Content.Find.Execute("45+6=51", $true, $true, $false, $false, $false, $true, 1, $false, "39+31=70", 2) | Out-Null
$d.Content.Find.Execute("90-4=86", $true, $true, $false, $false, $false, $true, 1, $false, "38+25=63", 2) | Out-Null
$d.Content.Find.Execute("44-26=18", $true, $true, $false, $false, $false, $true, 1, $false, "80-35=45", 2) | Out-Null
$d.Content.Find.Execute("93-59=34", $true, $true, $false, $false, $false, $true, 1, $false, "15+65=80", 2) | Out-Null
$d.Content.Find.Execute("26-13=13", $true, $true, $false, $false, $false, $true, 1, $false, "58-41=17", 2) | Out-Null
$d.Content.Find.Execute("84-54=30", $true, $true, $false, $false, $false, $true, 1, $false, "22+50=72", 2) | Out-Null
$d.Content.Find.Execute("97-10=87", $true, $true, $false, $false, $false, $true, 1, $false, "25+53=78", 2) | Out-Null
$d.Content.Find.Execute("21+17=38", $true, $true, $false, $false, $false, $true, 1, $false, "64-18=46", 2) | Out-Null
$d.Content.Find.Execute("56-55=1", $true, $true, $false, $false, $false, $true, 1, $false, "25-1=24", 2) | Out-Null
$d.Content.Find.Execute("62-17=45", $true, $true, $false, $false, $false, $true, 1, $false, "29-21=8", 2) | Out-Null
$d.Content.Find.Execute("50-43=7", $true, $true, $false, $false, $false, $true, 1, $false, "5+40=45", 2) | Out-Null
$d.Content.Find.Execute("46-28=18", $true, $true, $false, $false, $false, $true, 1, $false, "73-64=9", 2) | Out-Null
$d.Content.Find.Execute("1+69=70", $true, $true, $false, $false, $false, $true, 1, $false, "99-13=86", 2) | Out-Null
$d.Content.Find.Execute("76-28=48", $true, $true, $false, $false, $false, $true, 1, $false, "48+14=62", 2) | Out-Null
$d.Content.Find.Execute("16+58=74", $true, $true, $false, $false, $false, $true, 1, $false, "67-10=57", 2) | Out-Null
$d.Content.Find.Execute("84-37=47", $true, $true, $false, $false, $false, $true, 1, $false, "99-81=18", 2) | Out-Null
$d.Content.Find.Execute("10+14=24", $true, $true, $false, $false, $false, $true, 1, $false, "28+47=75", 2) | Out-Null
$d.Content.Find.Execute("91-69=22", $true, $true, $false, $false, $false, $true, 1, $false, "32+21=53", 2) | Out-Null
$d.Content.Find.Execute("64-23=41", $true, $true, $false, $false, $false, $true, 1, $false, "42-24=18", 2) | Out-Null
$d.Content.Find.Execute("16+80=96", $true, $true, $false, $false, $false, $true, 1, $false, "45-34=11", 2) | Out-Null
$d.Content.Find.Execute("94-25=69", $true, $true, $false, $false, $false, $true, 1, $false, "6+55=61", 2) | Out-Null
$d.Content.Find.Execute("31+21=52", $true, $true, $false, $false, $false, $true, 1, $false, "95-35=60", 2) | Out-Null
$d.Content.Find.Execute("17-6=11", $true, $true, $false, $false, $false, $true, 1, $false, "61-22=39", 2) | Out-Null
$d.Content.Find.Execute("0+68=68", $true, $true, $false, $false, $false, $true, 1, $false, "90-9=81", 2) | Out-Null
$d.Content.Find.Execute("85-28=57", $true, $true, $false, $false, $false, $true, 1, $false, "37+20=57", 2) | Out-Null
$d.Content.Find.Execute("8+8=16", $true, $true, $false, $false, $false, $true, 1, $false, "79-38=41", 2) | Out-Null
$d.Content.Find.Execute("23+8=31", $true, $true, $false, $false, $false, $true, 1, $false, "82-9=73", 2) | Out-Null
$d.Content.Find.Execute("81+5=86", $true, $true, $false, $false, $false, $true, 1, $false, "90-18=72", 2) | Out-Null
$d.Content.Find.Execute("64-54=10", $true, $true, $false, $false, $false, $true, 1, $false, "68-12=56", 2) | Out-Null
$d.Content.Find.Execute("75+24=99", $true, $true, $false, $false, $false, $true, 1, $false, "54-33=21", 2) | Out-Null
$d.Content.Find.Execute("63-22=41", $true, $true, $false, $false, $false, $true, 1, $false, "79+0=79", 2) | Out-Null
$d.Content.Find.Execute("21+50=71", $true, $true, $false, $false, $false, $true, 1, $false, "17-7=10", 2) | Out-Null
$d.Content.Find.Execute("90-38=52", $true, $true, $false, $false, $false, $true, 1, $false, "5+20=25", 2) | Out-Null
$d.Content.Find.Execute("35-3=32", $true, $true, $false, $false, $false, $true, 1, $false, "57+11=68", 2) | Out-Null
$d.Content.Find.Execute("9+31=40", $true, $true, $false, $false, $false, $true, 1, $false, "94-55=39", 2) | Out-Null
$d.Content.Find.Execute("86-50=36", $true, $true, $false, $false, $false, $true, 1, $false, "39-28=11", 2) | Out-Null
$d.Content.Find.Execute("41+27=68", $true, $true, $false, $false, $false, $true, 1, $false, "34+11=45", 2) | Out-Null
$d.Content.Find.Execute("16+70=86", $true, $true, $false, $false, $false, $true, 1, $false, "54-47=7", 2) | Out-Null
$d.Content.Find.Execute("29-20=9", $true, $true, $false, $false, $false, $true, 1, $false, "75-43=32", 2) | Out-Null
$d.Content.Find.Execute("83-35=48", $true, $true, $false, $false, $false, $true, 1, $false, "83-39=44", 2) | Out-Null
$d.Content.Find.Execute("47-16=31", $true, $true, $false, $false, $false, $true, 1, $false, "77-1=76", 2) | Out-Null
$d.Content.Find.Execute("26+8=34", $true, $true, $false, $false, $false, $true, 1, $false, "53-0=53", 2) | Out-Null
$d.Content.Find.Execute("18+24=42", $true, $true, $false, $false, $false, $true, 1, $false, "22+24=46", 2) | Out-Null
$d.Content.Find.Execute("36+1=37", $true, $true, $false, $false, $false, $true, 1, $false, "50+17=67", 2) | Out-Null
$d.Content.Find.Execute("22-20=2", $true, $true, $false, $false, $false, $true, 1, $false, "83-43=40", 2) | Out-Null
$d.Content.Find.Execute("27+12=39", $true, $true, $false, $false, $false, $true, 1, $false, "60+26=86", 2) | Out-Null
$d.Content.Find.Execute("72-51=21", $true, $true, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("23+44=67", $true, $true, $false, $false, $false, $true, 1, $false, "89-47=42", 2) | Out-Null
$d.Content.Find.Execute("47+49=96", $true, $true, $false, $false, $false, $true, 1, $false, "13-5=8", 2) | Out-Null
$d.Content.Find.Execute("6+71=77", $true, $true, $false, $false, $false, $true, 1, $false, "1+82=83", 2) | Out-Null
$d.Content.Find.Execute("65-43=22", $true, $true, $false, $false, $false, $true, 1, $false, "32-7=25", 2) | Out-Null
$d.Content.Find.Execute("75+13=88", $true, $true, $false, $false, $false, $true, 1, $false, "78-74=4", 2) | Out-Null
$d.Content.Find.Execute("8+57=65", $true, $true, $false, $false, $false, $true, 1, $false, "3+79=82", 2) | Out-Null
$d.Content.Find.Execute("43-34=9", $true, $true, $false, $false, $false, $true, 1, $false, "56-11=45", 2) | Out-Null
$d.Content.Find.Execute("74+0=74", $true, $true, $false, $false, $false, $true, 1, $false, "74+1=75", 2) | Out-Null
